$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Kalvin Phillips -> Stefan Ortega
$ws.Range("A2").Value = 240
$ws.Range("B2").Value = "Stefan Ortega"
$ws.Range("C2").Value = "GK"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 80
$ws.Range("J2").Value = "Germany"
$ws.Range("K2").Value = 30
$ws.Range("O2").Value = "Med"
$ws.Range("Q2").Value = 2

# Row 3: Victor Lindelof -> Karim Adeyemi
$ws.Range("A3").Value = 514
$ws.Range("B3").Value = "Karim Adeyemi"
$ws.Range("C3").Value = "ST"
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 75
$ws.Range("H3").Value = "Bundesliga"
$ws.Range("I3").Value = "Dortmund"
$ws.Range("J3").Value = "Germany"
$ws.Range("K3").Value = 21
$ws.Range("L3").Value = "Rare"
$ws.Range("M3").Value = "Left"
$ws.Range("N3").Value = "High"
$ws.Range("O3").Value = "Med"
$ws.Range("P3").Value = 650
$ws.Range("Q3").Value = 3

# Row 4: Sven Botman -> Niklas Stark
$ws.Range("A4").Value = 812
$ws.Range("B4").Value = "Niklas Stark"
$ws.Range("F4").Value = 75
$ws.Range("H4").Value = "Bundesliga"
$ws.Range("I4").Value = "Werder Bremen"
$ws.Range("J4").Value = "Germany"
$ws.Range("K4").Value = 27
$ws.Range("L4").Value = "Rare"
$ws.Range("M4").Value = "Right"
$ws.Range("O4").Value = "Med"
$ws.Range("Q4").Value = 3

# Row 5: Renan Lodi -> Marc-Oliver Kempf
$ws.Range("A5").Value = 2016
$ws.Range("B5").Value = "Marc-Oliver Kempf"
$ws.Range("C5").Value = "CB"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 75
$ws.Range("H5").Value = "Bundesliga"
$ws.Range("I5").Value = "Hertha Berlin"
$ws.Range("J5").Value = "Germany"
$ws.Range("K5").Value = 27
$ws.Range("L5").Value = "Rare"
$ws.Range("M5").Value = "Left"
$ws.Range("N5").Value = "High"
$ws.Range("O5").Value = "Med"
$ws.Range("P5").Value = 650
$ws.Range("Q5").Value = 3

# Row 6: Sergi Roberto -> Ridle Baku
$ws.Range("A6").Value = 2101
$ws.Range("B6").Value = "Ridle Baku"
$ws.Range("C6").Value = "RM"
$ws.Range("E6").Value = 4
$ws.Range("H6").Value = "Bundesliga"
$ws.Range("I6").Value = "VfL Wolfsburg"
$ws.Range("J6").Value = "Germany"
$ws.Range("K6").Value = 24
$ws.Range("L6").Value = "Rare"
$ws.Range("N6").Value = "High"
$ws.Range("Q6").Value = 3

# Row 7: Kangin Lee -> David Neres
$ws.Range("A7").Value = 2599
$ws.Range("B7").Value = "David Neres"
$ws.Range("C7").Value = "RW"
$ws.Range("D7").Value = 5
$ws.Range("F7").Value = 79
$ws.Range("G7").Value = "gold"
$ws.Range("H7").Value = "Liga NOS (POR 1)"
$ws.Range("I7").Value = "SL Benfica"
$ws.Range("J7").Value = "Brazil"
$ws.Range("K7").Value = 25
$ws.Range("O7").Value = "Med"

# Row 8: Oscar Trejo -> Lucas Verissimo
$ws.Range("A8").Value = 2601
$ws.Range("B8").Value = "Lucas Veríssimo"
$ws.Range("C8").Value = "CB"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 2
$ws.Range("H8").Value = "Liga NOS (POR 1)"
$ws.Range("I8").Value = "SL Benfica"
$ws.Range("J8").Value = "Brazil"
$ws.Range("K8").Value = 27
$ws.Range("O8").Value = "High"
$ws.Range("Q8").Value = 2

# Row 9: Youssef En-Nesyri -> Petr Sevcik
$ws.Range("A9").Value = 3003
$ws.Range("B9").Value = "Petr Ševčík"
$ws.Range("C9").Value = "CM"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 75
$ws.Range("H9").Value = "Česká Liga (CZE 1)"
$ws.Range("I9").Value = "Slavia Praha"
$ws.Range("J9").Value = "Czech Republic"
$ws.Range("K9").Value = 28
$ws.Range("L9").Value = "Rare"
$ws.Range("M9").Value = "Right"
$ws.Range("P9").Value = 650
# O9 unchanged ("High")

# Row 10: Mauro Arambarri -> Lukas Provod
$ws.Range("A10").Value = 3007
$ws.Range("B10").Value = "Lukas Provod"
$ws.Range("C10").Value = "CM"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 75
$ws.Range("H10").Value = "Česká Liga (CZE 1)"
$ws.Range("I10").Value = "Slavia Praha"
$ws.Range("J10").Value = "Czech Republic"
$ws.Range("K10").Value = 26
$ws.Range("L10").Value = "Rare"
$ws.Range("M10").Value = "Left"
$ws.Range("N10").Value = "High"
$ws.Range("O10").Value = "High"
$ws.Range("P10").Value = 650

# Row 11: Lucas Robertone -> Salem Al Dawsari
$ws.Range("A11").Value = 4780
$ws.Range("B11").Value = "Salem Al Dawsari"
$ws.Range("C11").Value = "LW"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 77
$ws.Range("H11").Value = "MBS Pro League (SAU 1)"
$ws.Range("I11").Value = "Al Hilal"
$ws.Range("J11").Value = "Saudi Arabia"
$ws.Range("K11").Value = 31
$ws.Range("N11").Value = "High"
# O11 unchanged ("High")
$ws.Range("Q11").Value = 0

# Row 12: Ledesma -> Ruben Vargas
$ws.Range("A12").Value = 9424
$ws.Range("B12").Value = "Ruben Vargas"
$ws.Range("C12").Value = "LM"
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 75
$ws.Range("H12").Value = "Bundesliga"
$ws.Range("I12").Value = " FC Augsburg"
$ws.Range("J12").Value = "Switzerland"
$ws.Range("K12").Value = 24
$ws.Range("L12").Value = "Rare"
$ws.Range("N12").Value = "High"
$ws.Range("O12").Value = "Med"
$ws.Range("P12").Value = 650
$ws.Range("Q12").Value = 2

$wb.Save()
